$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 23; I = 'ba'; J = 'Appreciation' }
    @{ Row = 39; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 44; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 50; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 55; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 57; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 74; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 76; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 98; I = '%'; J = 'Uninterpretable' }
    @{ Row = 112; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 113; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 117; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 120; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 121; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 124; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 127; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 138; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 143; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 145; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 164; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 174; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 175; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 176; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 180; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 195; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 196; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 197; I = 'ba'; J = 'Appreciation' }
    @{ Row = 210; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 228; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 229; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 238; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 240; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 243; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 250; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 251; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 253; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 255; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 258; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 259; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 273; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 274; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 275; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 277; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 309; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 327; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 332; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 354; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 359; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 376; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 382; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 388; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 389; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 391; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 393; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 402; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 403; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 409; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 413; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 416; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 420; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 428; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 448; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 460; I = '%'; J = 'Uninterpretable' }
    @{ Row = 465; I = '%'; J = 'Uninterpretable' }
    @{ Row = 469; I = '%'; J = 'Uninterpretable' }
    @{ Row = 474; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 475; I = 'ba'; J = 'Appreciation' }
    @{ Row = 478; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 488; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 492; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 496; I = 'sd'; J = 'Statement-non-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

Write-Output ("Updated {0} rows" -f $updates.Count)